# Generate Report for Handback
#
# This applies the "handback" update to localization-status.xlsx:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#     on every sheet that tracks it (Overview, zh-cn, de-de).
#   - The zh-cn and de-de detail sheets gain "Latest Target File" / "Latest
#     Handback File" hyperlink cells (columns F/G) for both data rows, linking
#     to the same package file / translated xlf file already linked from
#     columns A/D.
#   - The "Latest Handback DateTime" column (H) is stamped with the real
#     handback time per language (they used to share the 0001-01-01 sentinel).

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$mdFile    = "485d82c9-147c-4764-898d-eebdbf965e87.md"
$mdUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/6d9ed12ad85f442657b9c565218db780b8218e7e/e2e/$mdFile"

$zhXlf     = "485d82c9-147c-4764-898d-eebdbf965e87.61d42102b5c125f2746d19099ff0675d211731f5.zh-cn.xlf"
$zhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d96f6342e1f82239d2d5c39b552de5fc7ddfc5e2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf"

$deXlf     = "485d82c9-147c-4764-898d-eebdbf965e87.61d42102b5c125f2746d19099ff0675d211731f5.de-de.xlf"
$deXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/af076e0a521c6f15f57d623dbe954a4c7eb8ecae/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf"

$zhHandbackTime = "2016-03-12 16:50:53"
$deHandbackTime = "2016-03-12 16:50:59"

$hyperlinkRgb = 15570276  # matches the workbook's custom HyperLink font color FF6495ED

function Style-AsHyperlink($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hyperlinkRgb
}

# ---- Overview sheet: refresh the status text (cascades to every cell that
# shares the "Ready for handoff" string) ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus

$zh.Hyperlinks.Add($zh.Range("F2"), $mdUrl, "", "", $mdFile)
Style-AsHyperlink $zh.Range("F2")
$zh.Hyperlinks.Add($zh.Range("G2"), $zhXlfUrl, "", "", $zhXlf)
Style-AsHyperlink $zh.Range("G2")

$zh.Hyperlinks.Add($zh.Range("F3"), $mdUrl, "", "", $mdFile)
Style-AsHyperlink $zh.Range("F3")
$zh.Hyperlinks.Add($zh.Range("G3"), $zhXlfUrl, "", "", $zhXlf)
Style-AsHyperlink $zh.Range("G3")

$zh.Range("H2").Value = $zhHandbackTime
$zh.Range("H3").Value = $zhHandbackTime

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

$de.Hyperlinks.Add($de.Range("F2"), $mdUrl, "", "", $mdFile)
Style-AsHyperlink $de.Range("F2")
$de.Hyperlinks.Add($de.Range("G2"), $deXlfUrl, "", "", $deXlf)
Style-AsHyperlink $de.Range("G2")

$de.Hyperlinks.Add($de.Range("F3"), $mdUrl, "", "", $mdFile)
Style-AsHyperlink $de.Range("F3")
$de.Hyperlinks.Add($de.Range("G3"), $deXlfUrl, "", "", $deXlf)
Style-AsHyperlink $de.Range("G3")

$de.Range("H2").Value = $deHandbackTime
$de.Range("H3").Value = $deHandbackTime
